$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1) Replace the " : NoSettings" suffix with " : test" in the row
#    labels (column A) of the "BAU Emissions" sheet.
# ------------------------------------------------------------------
$wsBAU = $wb.Worksheets.Item("BAU Emissions")
for ($r = 1; $r -le 300; $r++) {
    $cell = $wsBAU.Cells.Item($r, 1)
    $val = $cell.Value()
    if ($val -ne $null -and $val -like "*: NoSettings") {
        $cell.Value = $val.Replace(" : NoSettings", " : test")
    }
}

# ------------------------------------------------------------------
# 2) Update the forecast values for row 94 (natural gas if, iron and
#    steel 241) columns M:AE.
# ------------------------------------------------------------------
$wsBAU.Range("M94").Value = 1001080
$wsBAU.Range("N94").Value = 2002150
$wsBAU.Range("O94").Value = 3003230
$wsBAU.Range("P94").Value = 4004300
$wsBAU.Range("Q94").Value = 5005380
$wsBAU.Range("R94").Value = 5005380
$wsBAU.Range("S94").Value = 5005380
$wsBAU.Range("T94").Value = 5005380
$wsBAU.Range("U94").Value = 5005380
$wsBAU.Range("V94").Value = 5005380
$wsBAU.Range("W94").Value = 5005380
$wsBAU.Range("X94").Value = 5005380
$wsBAU.Range("Y94").Value = 5005380
$wsBAU.Range("Z94").Value = 5005380
$wsBAU.Range("AA94").Value = 5005380
$wsBAU.Range("AB94").Value = 5005380
$wsBAU.Range("AC94").Value = 5005380
$wsBAU.Range("AD94").Value = 5005380
$wsBAU.Range("AE94").Value = 5005380

# ------------------------------------------------------------------
# 3) Update the "About" sheet's printed date (C1) and make it the
#    active/selected sheet.
# ------------------------------------------------------------------
$wsAbout = $wb.Worksheets.Item("About")
$wsAbout.Range("C1").Value = 45387

# ------------------------------------------------------------------
# 4) Update the "BAU Emissions" sheet's current selection/scroll
#    position.
# ------------------------------------------------------------------
$wsBAU.Range("A30:AE280").Select()

# ------------------------------------------------------------------
# 5) Activate "About" last so it becomes the workbook's active tab
#    (moving the selection away from "Current and Planned Capacity").
# ------------------------------------------------------------------
$wsAbout.Activate()
